# Preparing repository for submission
#
# The CST storage block used two "Link" rows (heat storage link in / out)
# plus a "Store" row to model the molten-salt storage. Simplify this to a
# single PyPSA "StorageUnit" row and drop the now-unused heat-storage link
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "heat storage link in" / "heat storage link out" rows
# (old rows 53 and 54). This shifts the old "Store" row (55) up to become
# the new row 53, and everything below moves up by two as well.
$ws.Rows("53:54").Delete()

# The row that used to describe the molten-salt "Store" now becomes a
# PyPSA "StorageUnit" component, and it connects directly to the "cst_out"
# bus instead of the removed "heat_storage" bus.
$ws.Range("A53").Value = "StorageUnit"
$ws.Range("D53").Value = "cst_out"

# Restore the user's on-screen selection/view state.
$ws.Range("C56").Select() | Out-Null
